$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'87.793.24"
$ws.Range("E2").Value = "'  -2.34%  "
$ws.Range("D3").Value = "'3.050.89"
$ws.Range("E3").Value = "'  -5.13%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'209.02"
$ws.Range("E5").Value = "'  -4.51%  "
$ws.Range("D6").Value = "'615.81"
$ws.Range("E6").Value = "'  -2.42%  "
$ws.Range("D7").Value = "'0.367"
$ws.Range("E7").Value = "'  -7.24%  "
$ws.Range("D8").Value = "'0.799"
$ws.Range("E8").Value = "'  +14.33%  "
$ws.Range("D10").Value = "'3.049.71"
$ws.Range("E10").Value = "'  -5.18%  "
$ws.Range("D11").Value = "'0.597"
$ws.Range("E11").Value = "'  +2.32%  "
$ws.Range("D12").Value = "'0.178"
$ws.Range("E12").Value = "'  -1.16%  "
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "'  -10.73%  "
$ws.Range("D14").Value = "'5.26"
$ws.Range("E14").Value = "'  -3.57%  "
$ws.Range("D15").Value = "'87.613.36"
$ws.Range("E15").Value = "'  -2.51%  "
$ws.Range("D16").Value = "'3.616.50"
$ws.Range("E16").Value = "'  -5.14%  "
$ws.Range("D17").Value = "'31.73"
$ws.Range("E17").Value = "'  -6.91%  "
$ws.Range("D18").Value = "'3.052.77"
$ws.Range("E18").Value = "'  -5.84%  "
$ws.Range("D19").Value = "'3.23"
$ws.Range("E19").Value = "'  -6.29%  "
$ws.Range("D20").Value = "'0.0000198"
$ws.Range("E20").Value = "'  -15.69%  "
$ws.Range("D21").Value = "'13.16"
$ws.Range("E21").Value = "'  -3.90%  "
$ws.Range("D22").Value = "'418.23"
$ws.Range("E22").Value = "'  -5.27%  "
$ws.Range("D23").Value = "'8.10"
$ws.Range("E23").Value = "'  -7.57%  "
$ws.Range("D24").Value = "'4.85"
$ws.Range("E24").Value = "'  -5.99%  "
$ws.Range("D25").Value = "'5.45"
$ws.Range("E25").Value = "'  +2.13%  "
$ws.Range("D26").Value = "'11.66"
$ws.Range("E26").Value = "'  -4.83%  "
$ws.Range("D27").Value = "'81.93"
$ws.Range("E27").Value = "'  -2.12%  "
$ws.Range("B28").Value = "'WrappedeETH"
$ws.Range("C28").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'3.228.09"
$ws.Range("E28").Value = "'  -6.52%  "
$ws.Range("B29").Value = "'Dai"
$ws.Range("C29").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("B30").Value = "'Binance-PegBSC-USD"
$ws.Range("C30").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.09"
$ws.Range("E30").Value = "'  +9.05%  "
$ws.Range("B31").Value = "'Cronos"
$ws.Range("C31").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "'0.168"
$ws.Range("E31").Value = "'  +4.47%  "
$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.00"
$ws.Range("E32").Value = "'  -7.01%  "
$ws.Range("B33").Value = "'Bittensor"
$ws.Range("C33").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'505.68"
$ws.Range("E33").Value = "'  -8.47%  "
$ws.Range("B34").Value = "'dogwifhat"
$ws.Range("C34").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "'3.53"
$ws.Range("E34").Value = "'  -15.99%  "
$ws.Range("B35").Value = "'RenderToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").Value = "'6.69"
$ws.Range("E35").Value = "'  -6.74%  "
$ws.Range("B36").Value = "'PancakeSwap"
$ws.Range("C36").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "'1.78"
$ws.Range("E36").Value = "'  -7.90%  "
$ws.Range("B37").Value = "'Fetch.AI"
$ws.Range("C37").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'1.24"
$ws.Range("E37").Value = "'  -7.60%  "
$ws.Range("B38").Value = "'EthereumClassic"
$ws.Range("C38").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'22.16"
$ws.Range("E38").Value = "'  -2.02%  "
$ws.Range("B39").Value = "'Kaspa"
$ws.Range("C39").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.131"
$ws.Range("E39").Value = "'  +0.42%  "
$ws.Range("B40").Value = "'WhiteBITCoin"
$ws.Range("C40").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'22.21"
$ws.Range("E40").Value = "'  -0.88%  "
$ws.Range("B41").Value = "'FirstDigitalUSD"
$ws.Range("C41").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E41").Value = "'  +0.34%  "
$ws.Range("B42").Value = "'USDe"
$ws.Range("C42").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "'  +0.02%  "
$ws.Range("B43").Value = "'PolygonEcosystemToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.358"
$ws.Range("E43").Value = "'  -5.89%  "
$ws.Range("B44").Value = "'Monero"
$ws.Range("C44").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'147.18"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("E45").Value = "'  +5.42%  "
$ws.Range("B46").Value = "'Stacks"
$ws.Range("C46").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.78"
$ws.Range("E46").Value = "'  -9.08%  "
$ws.Range("B47").Value = "'OKB"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'43.37"
$ws.Range("E47").Value = "'  -1.06%  "
$ws.Range("B48").Value = "'Hedera"
$ws.Range("C48").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0670"
$ws.Range("E48").Value = "'  +9.76%  "
$ws.Range("B49").Value = "'Mantle"
$ws.Range("C49").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.698"
$ws.Range("E49").Value = "'  -10.62%  "
$ws.Range("D50").Value = "'154.87"
$ws.Range("E50").Value = "'  -11.67%  "
$ws.Range("B51").Value = "'ImmutableX"
$ws.Range("C51").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").Value = "'1.17"
$ws.Range("E51").Value = "'  -7.41%  "
